$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.3464964993005633
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 1.336873824401267

# Row 3
$ws.Range("B3").Value = 0.1554434735375247
$ws.Range("C3").Value = 0.004309184025731883
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 0.8125451466554947

# Row 4
$ws.Range("B4").Value = 0.3464964993005633
$ws.Range("C4").Value = 0.3375848360084654
$ws.Range("D4").Value = 0.1529057820181812
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 1.336873824401267

# Row 5
$ws.Range("B5").Value = 3.182878228561681
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 0.1529057820181812
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("G5").Value = 5.488907176552729
